$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("K2").Value = "30-05-2024"
$ws1.Range("N2").Value = "30-05-2024 02:21:31 PM"
$ws1.Range("AG2").Value = "ET464"
$ws1.Range("AK2").Value = "'3"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("K2").Value = "30-05-2024"
$ws2.Range("N2").Value = "30-05-2024 02:27:31 PM"
$ws2.Range("AG2").Value = "ET465"

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("K2").Value = "30-05-2024"
$ws3.Range("N2").Value = "30-05-2024 02:27:31 PM"
$ws3.Range("AG2").Value = "ET466"

$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("K2").Value = "30-05-2024"
$ws4.Range("N2").Value = "30-05-2024 02:27:31 PM"
$ws4.Range("AG2").Value = "ET466"
